$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New goal values for the remaining Goal1 column (B)
$ws.Range("B2").Value = 232
$ws.Range("B3").Value = 148

# Remove the other goal columns (Goal2/Goal3/Goal4 data) - clear contents only,
# preserving any cell formatting (matches the observed quotePrefix style cells)
$ws.Range("C1:E8").ClearContents()
$ws.Range("C9:E9").ClearContents()

# Update the active selection to B4
$ws.Range("B4").Select()
